$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AD1:AF1,
# matching the style of the existing header cells (e.g. AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the team record (Wins/Losses/Ties) for every player row.
for ($row = 2; $row -le 42; $row++) {
    $ws.Cells.Item($row, 30).Value = 75
    $ws.Cells.Item($row, 31).Value = 86
    $ws.Cells.Item($row, 32).Value = 0
}
